$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

$lines = @(
    "Today:",
    "Recursion (Ch 5 think python) (Workbook 008)",
    "Some large exercises mixed in. ",
    "Overall Context:",
    "We’ve gone over most of the basic tools we need to ”do stuff” ",
    "We’ll circle back around to some of those concepts, with a more full set of tools, and try to build things that are more complex bit by bit. ",
    "I think it is much easier, in the long run, to learn through larger scoped tasks. ",
    "We have most of the key building blocks of programs, now need to practice. "
)

# IndentLevel 1 == outline level 0 (no <a:pPr lvl>), IndentLevel 2 == outline level 1 (<a:pPr lvl="1"/>)
$levels = @(1, 2, 2, 1, 2, 2, 2, 2)

$tr.Text = [string]::Join([char]13, $lines)

for ($i = 1; $i -le $lines.Count; $i++) {
    if ($levels[$i - 1] -eq 2) {
        $tr.Paragraphs($i, 1).IndentLevel = 2
    }
}
